$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1347866666666667
$ws.Range("H2").Value = 0.40436
$ws.Range("I2").Value = 0.03419045085634245
$ws.Range("J2").Value = 0.03419045085634244
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.05529753293777778
$ws.Range("R2").Value = 0.49767779644
$ws.Range("S2").Value = 0.0001196535812947071
$ws.Range("T2").Value = 0.0001196535812947071

$ws.Range("G3").Value = 0.1347866666666667
$ws.Range("H3").Value = 0.40436
$ws.Range("I3").Value = 0.03419045085634245
$ws.Range("J3").Value = 0.03419045085634244
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 13.73451507809333
$ws.Range("R3").Value = 123.61063570284
$ws.Range("S3").Value = 0.02971893734010158
$ws.Range("T3").Value = 0.02971893734010157

$ws.Range("G4").Value = 0.1347866666666667
$ws.Range("H4").Value = 0.40436
$ws.Range("I4").Value = 0.03419045085634245
$ws.Range("J4").Value = 0.03419045085634244
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 2.011198624306667
$ws.Range("R4").Value = 18.10078761876
$ws.Range("S4").Value = 0.004351859934946161
$ws.Range("T4").Value = 0.00435185993494616

$ws.Range("I5").Value = 0.3318597741685039
$ws.Range("J5").Value = 0.3318597741685039
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 0.5367295935906667
$ws.Range("R5").Value = 4.830566342316001
$ws.Range("S5").Value = 0.001161383060836363
$ws.Range("T5").Value = 0.001161383060836362

$ws.Range("I6").Value = 0.3318597741685039
$ws.Range("J6").Value = 0.3318597741685039
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.2884583147397864
$ws.Range("T6").Value = 0.2884583147397863

$ws.Range("I7").Value = 0.3318597741685039
$ws.Range("J7").Value = 0.3318597741685039
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("S7").Value = 0.04224007636788117
$ws.Range("T7").Value = 0.04224007636788116

$ws.Range("I8").Value = 0.6339497749751537
$ws.Range("J8").Value = 0.6339497749751537
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 1.025311386207778
$ws.Range("R8").Value = 9.227802475870002
$ws.Range("S8").Value = 0.002218583231191278
$ws.Range("T8").Value = 0.002218583231191278

$ws.Range("I9").Value = 0.6339497749751537
$ws.Range("J9").Value = 0.6339497749751537
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("S9").Value = 0.5510402222661287
$ws.Range("T9").Value = 0.5510402222661286

$ws.Range("I10").Value = 0.6339497749751537
$ws.Range("J10").Value = 0.6339497749751537
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("Q10").Value = 37.29108225763666
$ws.Range("R10").Value = 335.61974031873
$ws.Range("S10").Value = 0.08069096947783383
$ws.Range("T10").Value = 0.08069096947783383
